$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue {
    param($addr, $val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "244.62"
Set-TextValue "D3" "21.87"
Set-TextValue "D4" "5.398"
Set-TextValue "D5" "0.06040"
Set-TextValue "D6" "3.398"
Set-TextValue "D7" "0.8140"
Set-TextValue "D10" "0.07472"
Set-TextValue "D11" "0.03362"
Set-TextValue "D12" "0.03049"
Set-TextValue "D13" "4.008"
Set-TextValue "D14" "0.09403"
Set-TextValue "D15" "0.001595"
Set-TextValue "D16" "0.04810"
Set-TextValue "D17" "0.0005942"
Set-TextValue "D18" "0.005394"
Set-TextValue "D19" "0.004157"
Set-TextValue "D20" "0.0009905"
Set-TextValue "D22" "3.653"
Set-TextValue "D23" "6.434"
Set-TextValue "D27" "0.0002901"
Set-TextValue "D40" "0.03983"
Set-TextValue "D41" "0.006402"
Set-TextValue "D42" "0.1076"
Set-TextValue "D44" "0.006373"
Set-TextValue "D45" "0.00005230"
Set-TextValue "D48" "0.002527"
Set-TextValue "D49" "0.00002101"
